$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Mark the "Database" section tasks (rows 13-16) as completed ---
# Column F holds a boolean "Status" checkbox; flip it from FALSE to TRUE.
$ws.Range("F13").Value = $true
$ws.Range("F14").Value = $true
$ws.Range("F15").Value = $true
$ws.Range("F16").Value = $true

# --- Highlight all of the completed task rows (rows 2-16, columns A-E) ---
# with a green fill, thin border, and vertical-center + wrap-text alignment.
$rng = $ws.Range("A2:E16")
$rng.Interior.Color = 5296274   # RGB(146, 208, 80) -> 0x92D050, stored BGR for COM
$rng.Borders.LineStyle = 1
$rng.VerticalAlignment = -4108  # xlCenter
$rng.WrapText = $true

# --- Update the sheet view / selection to match where the user was working ---
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D14").Select()
